# sao luu du lieu ngay 02/06/2025
# Replaces the single declared-guest row with an updated guest record and
# appends a full batch of re-declarations (rows 2-12), mirroring how the
# source "DU_LIEU_KHAI_BAO" log accumulates one row per saved CCCD photo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- shared, per-record field values (same guest, same ID photos) ----
$soGiayTo   = 91203001544
$soCmndCu   = 371987447
$hoVaTen    = "Nguyễn Văn Hoài"
$gioiTinh   = "Nam"
$ngaySinh   = "23/02/2003"
$noiThuongTru = "Khu Phố 4, An Thới, TP. Phú Quốc, Kiên Giang"
$ngayCap    = 44381
$loaiGiayTo = "CCCD"
$roomOld    = "Phòng 3 nhà cũ"
$roomNew    = "Phòng 4 nhà cũ"
$basePath   = "d:\QUAN LY NHA NGHI\APP_KHAI_BAO_LUU_TRU_2\Anh_CCCD_da_khai_bao\"

# Per-row variable data: room name, "Thoi gian ghi" timestamp serial, and the
# timestamped filename suffix shared by the front/back photo pair.
$rows = @(
  @{ room = $roomOld; t = 45694.479907407411; suf = "20250602_113104" },
  @{ room = $roomOld; t = 45694.481724537036; suf = "20250602_113340" },
  @{ room = $roomOld; t = 45694.483067129629; suf = "20250602_113537" },
  @{ room = $roomOld; t = 45694.484884259262; suf = "20250602_113814" },
  @{ room = $roomOld; t = 45694.485289351855; suf = "20250602_113849" },
  @{ room = $roomOld; t = 45694.486400462964; suf = "20250602_114025" },
  @{ room = $roomNew; t = 45694.486956018518; suf = "20250602_114113" },
  @{ room = $roomOld; t = 45694.488530092596; suf = "20250602_114329" },
  @{ room = $roomOld; t = 45694.489305555559; suf = "20250602_114436" },
  @{ room = $roomOld; t = 45694.49255787037;  suf = "20250602_114917" },
  @{ room = $roomNew; t = 45694.493321759262; suf = "20250602_115023" }
)

$firstRow = 2
$lastRow  = $firstRow + $rows.Count - 1   # 2 .. 12

# ---- 1. overwrite row 2 in place with the updated guest data ----
$ws.Cells.Item($firstRow, 1).Value = $soGiayTo
$ws.Cells.Item($firstRow, 2).Value = $soCmndCu
$ws.Cells.Item($firstRow, 3).Value = $hoVaTen
$ws.Cells.Item($firstRow, 4).Value = $gioiTinh
$ws.Cells.Item($firstRow, 5).Value = $ngaySinh
$ws.Cells.Item($firstRow, 6).Value = $noiThuongTru
$ws.Cells.Item($firstRow, 7).Value = $ngayCap
$ws.Cells.Item($firstRow, 7).NumberFormat = "mm-dd-yy"
$ws.Cells.Item($firstRow, 8).Value = $loaiGiayTo

# ---- 2. propagate row 2's full formatting down through row 12 ----
# (keeps fonts / alignment / borders / hyperlink style identical to the
# original single-row template, including the now-dated G column format)
$srcRow = $ws.Range("A" + $firstRow + ":L" + $firstRow)
$srcRow.Copy()
$dstRows = $ws.Range("A" + ($firstRow + 1) + ":L" + $lastRow)
$dstRows.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---- 3. fill in every row's values / formulas ----
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $firstRow + $i
  $d = $rows[$i]

  $ws.Cells.Item($r, 1).Value = $soGiayTo
  $ws.Cells.Item($r, 2).Value = $soCmndCu
  $ws.Cells.Item($r, 3).Value = $hoVaTen
  $ws.Cells.Item($r, 4).Value = $gioiTinh
  $ws.Cells.Item($r, 5).Value = $ngaySinh
  $ws.Cells.Item($r, 6).Value = $noiThuongTru
  $ws.Cells.Item($r, 7).Value = $ngayCap
  $ws.Cells.Item($r, 8).Value = $loaiGiayTo
  $ws.Cells.Item($r, 9).Value = $d.room
  $ws.Cells.Item($r, 10).Value = $d.t

  $frontPath = $basePath + "mat_truoc_" + $d.suf + ".jpg"
  $backPath  = $basePath + "mat_sau_" + $d.suf + ".jpg"
  $ws.Cells.Item($r, 11).Formula = '=HYPERLINK("' + $frontPath + '", "Ảnh mặt trước")'
  $ws.Cells.Item($r, 12).Formula = '=HYPERLINK("' + $backPath + '", "Ảnh mặt sau")'
}

# ---- 4. cosmetic cleanup matching the refreshed table: resize columns ----
# whose best-fit width changed because of the new name / address text, and
# move the selection back onto the sheet's first cell.
$ws.Columns("C:C").AutoFit()
$ws.Columns("F:F").AutoFit()
$ws.Range("A1").Select()

Write-Host "Rewrote row $firstRow and appended rows $($firstRow+1)-$lastRow."
